$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing, so that numeric-looking
# strings (e.g. "1.008", "12.50") are preserved verbatim as text rather than being
# normalized into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.661.42'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '2.089.78'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').Value = '344.08'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -1.96%  '
$ws.Range('D8').Value = '0.4399'
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('D9').Value = '0.09269'
$ws.Range('E9').Value = '  +2.90%  '
$ws.Range('D10').Value = '52.03'
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('D11').Value = '1.177'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '25.19'
$ws.Range('E12').Value = '  +2.78%  '
$ws.Range('D13').Value = '2.094.48'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '6.764'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '8.169'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '100.11'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').Value = '0.00001158'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').Value = '21.02'
$ws.Range('E19').Value = '  +8.64%  '
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '6.200'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').Value = '29.710.16'
$ws.Range('E23').Value = '  -3.23%  '
$ws.Range('D24').Value = '12.71'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').Value = '2.309'
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('D26').Value = '2.333.31'
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('D27').Value = '21.90'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('D28').Value = '162.92'
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').Value = '2.524'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '132.76'
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('D31').Value = '1.143'
$ws.Range('E31').Value = '  -3.79%  '
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').Value = '1.631'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').Value = '6.192'
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('D35').Value = '3.972'
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').Value = '6.093'
$ws.Range('E36').Value = '  +3.17%  '
$ws.Range('D37').Value = '10.33'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').Value = '0.02572'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = '0.06715'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '12.50'
$ws.Range('E40').Value = '  -0.65%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.2248'
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('D42').Value = '0.6857'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = '1.297'
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('D44').Value = '0.6647'
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').Value = '14.24'
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').Value = '2.336'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('D47').Value = '3.623'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000336'
$ws.Range('E49').Value = '  -7.15%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '81.62'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').Value = '1.168'
$ws.Range('E51').Value = '  -2.28%  '

# Restore the default "Normal" style on column D so no stray explicit cell style
# (number format) is left behind compared to the original workbook.
$ws.Range("D2:D51").Style = "Normal"

